$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 672.8570999999999
$ws.Range("J2").Value = 600
$ws.Range("L2").Value = 600
$ws.Range("N2").Value = -826
$ws.Range("H38").Value = 4167.6665
$ws.Range("H64").Value = 3200
$ws.Range("J64").Value = 3200
$ws.Range("L64").Value = 3200
$ws.Range("N64").Value = -3696
$ws.Range("H67").Value = 3200
$ws.Range("J67").Value = 3200
$ws.Range("L67").Value = 3200
$ws.Range("N67").Value = -4916
$ws.Range("H115").Value = 5541.7144
$ws.Range("I115").Value = 5541.7144
$ws.Range("K115").Value = 16625.1432
$ws.Range("M115").Value = -15058.1432
$ws.Range("H138").Value = 60098
$ws.Range("J138").Value = 69999
$ws.Range("L138").Value = 209997
$ws.Range("N138").Value = -220277

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13729.692
$ws.Range("I32").Value = 5610
$ws.Range("J32").Value = 24237.53
$ws.Range("K32").Value = 5610
$ws.Range("L32").Value = 24237.53
$ws.Range("M32").Value = -5323
$ws.Range("N32").Value = -24811.53
$ws.Range("H61").Value = 1613.909
$ws.Range("I61").Value = 1650.4445
$ws.Range("J61").Value = 1449.5
$ws.Range("K61").Value = 1650.4445
$ws.Range("L61").Value = 1449.5
$ws.Range("M61").Value = -1438.4445
$ws.Range("N61").Value = -1873.5
$ws.Range("H63").Value = 7367.2
$ws.Range("I63").Value = 6358.5713
$ws.Range("K63").Value = 6358.5713
$ws.Range("M63").Value = -5672.5713
$ws.Range("H66").Value = 7367.2
$ws.Range("I66").Value = 6358.5713
$ws.Range("K66").Value = 31792.8565
$ws.Range("M66").Value = -28360.8565
$ws.Range("H74").Value = 3857
$ws.Range("I74").Value = 2010
$ws.Range("J74").Value = 4549.625
$ws.Range("K74").Value = 2010
$ws.Range("L74").Value = 4549.625
$ws.Range("M74").Value = -1136
$ws.Range("N74").Value = -6297.625
$ws.Range("H77").Value = 3857
$ws.Range("I77").Value = 2010
$ws.Range("J77").Value = 4549.625
$ws.Range("K77").Value = 10050
$ws.Range("L77").Value = 22748.125
$ws.Range("M77").Value = -5682
$ws.Range("N77").Value = -31484.125
$ws.Range("H136").Value = 1613.909
$ws.Range("I136").Value = 1650.4445
$ws.Range("J136").Value = 1449.5
$ws.Range("K136").Value = 4951.333500000001
$ws.Range("L136").Value = 4348.5
$ws.Range("M136").Value = -2401.333500000001
$ws.Range("N136").Value = -9448.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5170.1665
$ws.Range("I31").Value = 3799.6667
$ws.Range("J31").Value = 5855.4165
$ws.Range("K31").Value = 3799.6667
$ws.Range("L31").Value = 5855.4165
$ws.Range("M31").Value = -3504.6667
$ws.Range("N31").Value = -6445.4165
$ws.Range("H34").Value = 5170.1665
$ws.Range("I34").Value = 3799.6667
$ws.Range("J34").Value = 5855.4165
$ws.Range("K34").Value = 3799.6667
$ws.Range("L34").Value = 5855.4165
$ws.Range("M34").Value = -3597.6667
$ws.Range("N34").Value = -6259.4165
$ws.Range("H134").Value = 5306.091
$ws.Range("I134").Value = 4302.125
$ws.Range("K134").Value = 12906.375
$ws.Range("M134").Value = -10371.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 445.9643
$ws.Range("J5").Value = 833.3333
$ws.Range("L5").Value = 2499.9999
$ws.Range("N5").Value = -2723.9999
$ws.Range("H13").Value = 5400
$ws.Range("I13").Value = 800
$ws.Range("J13").Value = 10000
$ws.Range("K13").Value = 2400
$ws.Range("L13").Value = 30000
$ws.Range("M13").Value = -2232
$ws.Range("N13").Value = -30336
$ws.Range("H17").Value = 759.6667
$ws.Range("I17").Value = 340
$ws.Range("K17").Value = 1020
$ws.Range("M17").Value = -851
$ws.Range("H81").Value = 2318.5
$ws.Range("I81").Value = 999
$ws.Range("K81").Value = 2997
$ws.Range("M81").Value = -1874
$ws.Range("H84").Value = 2318.5
$ws.Range("I84").Value = 999
$ws.Range("K84").Value = 8991
$ws.Range("M84").Value = -3375
$ws.Range("H109").Value = 2000
$ws.Range("H113").Value = 2116.4
$ws.Range("J113").Value = 1808.125
$ws.Range("L113").Value = 5424.375
$ws.Range("N113").Value = -9764.375
$ws.Range("H122").Value = 1142.7142
$ws.Range("H135").Value = 445.9643
$ws.Range("J135").Value = 833.3333
$ws.Range("L135").Value = 7499.9997
$ws.Range("N135").Value = -12569.9997
$ws.Range("H137").Value = 5174.375
$ws.Range("I137").Value = 6000
$ws.Range("J137").Value = 4348.75
$ws.Range("K137").Value = 18000
$ws.Range("L137").Value = 13046.25
$ws.Range("M137").Value = -12900
$ws.Range("N137").Value = -23246.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1167.6666
$ws.Range("I3").Value = 3
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 113
$ws.Range("H80").Value = 4940.4
$ws.Range("I80").Value = 4039.6
$ws.Range("J80").Value = 5841.2
$ws.Range("K80").Value = 4039.6
$ws.Range("L80").Value = 5841.2
$ws.Range("M80").Value = -3041.6
$ws.Range("N80").Value = -7837.2
$ws.Range("H83").Value = 4940.4
$ws.Range("I83").Value = 4039.6
$ws.Range("J83").Value = 5841.2
$ws.Range("K83").Value = 20198
$ws.Range("L83").Value = 29206
$ws.Range("M83").Value = -15206
$ws.Range("N83").Value = -39190
$ws.Range("H100").Value = 42499.625
$ws.Range("J100").Value = 42499.625
$ws.Range("L100").Value = 42499.625
$ws.Range("N100").Value = -44663.625
$ws.Range("H107").Value = 1955.6666
$ws.Range("I107").Value = 9501
$ws.Range("J107").Value = 1012.5
$ws.Range("K107").Value = 9501
$ws.Range("L107").Value = 1012.5
$ws.Range("M107").Value = -7581
$ws.Range("N107").Value = -4852.5
$ws.Range("H126").Value = 4433.353
$ws.Range("I126").Value = 3379.5
$ws.Range("J126").Value = 5008.1816
$ws.Range("K126").Value = 10138.5
$ws.Range("L126").Value = 15024.5448
$ws.Range("M126").Value = -7668.5
$ws.Range("N126").Value = -19964.5448

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 9162.5
$ws.Range("I16").Value = 9162.5
$ws.Range("K16").Value = 9162.5
$ws.Range("M16").Value = -8992.5
$ws.Range("H40").Value = 3930.6843
$ws.Range("I40").Value = 4156.143
$ws.Range("J40").Value = 3299.4
$ws.Range("K40").Value = 4156.143
$ws.Range("L40").Value = 3299.4
$ws.Range("M40").Value = -4020.143
$ws.Range("N40").Value = -3571.4
$ws.Range("H61").Value = 5754
$ws.Range("I61").Value = 5963.222
$ws.Range("J61").Value = 5377.4
$ws.Range("K61").Value = 5963.222
$ws.Range("L61").Value = 5377.4
$ws.Range("M61").Value = -5761.222
$ws.Range("N61").Value = -5781.4
$ws.Range("H93").Value = 1708.5714
$ws.Range("I93").Value = 1212
$ws.Range("J93").Value = 2950
$ws.Range("K93").Value = 1212
$ws.Range("L93").Value = 2950
$ws.Range("M93").Value = 36
$ws.Range("N93").Value = -5446
$ws.Range("H100").Value = 4579.7144
$ws.Range("I100").Value = 1061.6
$ws.Range("J100").Value = 13375
$ws.Range("K100").Value = 1061.6
$ws.Range("L100").Value = 13375
$ws.Range("M100").Value = -520.5999999999999
$ws.Range("N100").Value = -14457
$ws.Range("H103").Value = 20000
$ws.Range("J103").Value = 20000
$ws.Range("L103").Value = 20000
$ws.Range("N103").Value = -22344
$ws.Range("H113").Value = 5754
$ws.Range("I113").Value = 5963.222
$ws.Range("J113").Value = 5377.4
$ws.Range("K113").Value = 5963.222
$ws.Range("L113").Value = 5377.4
$ws.Range("M113").Value = -3793.222
$ws.Range("N113").Value = -9717.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2243.7778
$ws.Range("I81").Value = 2243.7778
$ws.Range("K81").Value = 4487.5556
$ws.Range("M81").Value = -3426.5556
$ws.Range("H84").Value = 2243.7778
$ws.Range("I84").Value = 2243.7778
$ws.Range("K84").Value = 22437.778
$ws.Range("M84").Value = -17133.778
$ws.Range("H113").Value = 2112.25
$ws.Range("I113").Value = 1625.75
$ws.Range("J113").Value = 2598.75
$ws.Range("K113").Value = 4877.25
$ws.Range("L113").Value = 7796.25
$ws.Range("M113").Value = -2707.25
$ws.Range("N113").Value = -12136.25
$ws.Range("H126").Value = 41966.08
$ws.Range("I126").Value = 101624.8
$ws.Range("J126").Value = 2193.6
$ws.Range("K126").Value = 304874.4
$ws.Range("L126").Value = 6580.799999999999
$ws.Range("M126").Value = -302404.4
$ws.Range("N126").Value = -11520.8
$ws.Range("H132").Value = 1887.2069
$ws.Range("I132").Value = 1370.2632
$ws.Range("K132").Value = 4110.7896
$ws.Range("M132").Value = -1580.7896
$ws.Range("H136").Value = 114374.445
$ws.Range("I136").Value = 2910.8572
$ws.Range("K136").Value = 8732.571599999999
$ws.Range("M136").Value = -6182.571599999999
